$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a pure number (e.g. "706.40") need the
# cell pre-formatted as Text so Excel stores the literal text exactly as
# scraped, matching the original inline-string "Price" column cells.
$ws.Range("D2").Value = "70.610.48"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "3.792.52"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "706.40"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.20"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").Value = "3.791.70"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.10"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "4.435.94"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "3.840.86"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "70.664.60"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "492.79"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  -4.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.91"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").Value = "3.946.28"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -4.97%  "
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.06"
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "3.765.84"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.04"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("E43").Value = "  -4.16%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000318"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "164.57"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.87"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "420.25"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  -1.66%  "
